$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.271.09'
$ws.Range("E2").Value = '  +0.44%  '
$ws.Range("D3").Value = '1.662.13'
$ws.Range("E3").Value = '  +0.26%  '
$ws.Range("E4").Value = '  +0.77%  '
$ws.Range("D5").Value = '''218.37'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = '''0.5318'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.68%  '
$ws.Range("E7").Value = '  +0.71%  '
$ws.Range("D8").Value = '''0.2635'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.89%  '
$ws.Range("D9").Value = '''0.06358'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.06%  '
$ws.Range("D10").Value = '''20.53'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.21%  '
$ws.Range("D11").Value = '''0.07835'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.67%  '
$ws.Range("D12").Value = '''4.565'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.15%  '
$ws.Range("D13").Value = '1.667.51'
$ws.Range("E13").Value = '  +0.36%  '
$ws.Range("D14").Value = '1.889.98'
$ws.Range("E14").Value = '  +0.23%  '
$ws.Range("D15").Value = '''0.5528'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.43%  '
$ws.Range("D16").Value = '0.0₅8152'
$ws.Range("E16").Value = '  -1.17%  '
$ws.Range("D17").Value = '''65.60'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.10%  '
$ws.Range("E18").Value = '  +0.70%  '
$ws.Range("D19").Value = '''4.667'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.73%  '
$ws.Range("D20").Value = '''192.98'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.16%  '
$ws.Range("D21").Value = '''10.19'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.83%  '
$ws.Range("E22").Value = '  -0.38%  '
$ws.Range("E23").Value = '  +0.78%  '
$ws.Range("D24").Value = '''145.63'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.60%  '
$ws.Range("D25").Value = '''0.1220'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.22%  '
$ws.Range("D26").Value = '''7.188'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.35%  '
$ws.Range("E27").Value = '  -0.74%  '
$ws.Range("D28").Value = '''1.488'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.26%  '
$ws.Range("D29").Value = '''0.05890'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.71%  '
$ws.Range("E30").Value = '  -0.15%  '
$ws.Range("D31").Value = '''3.581'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.45%  '
$ws.Range("D32").Value = '''3.271'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.03%  '
$ws.Range("D33").Value = '''1.607'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.07%  '
$ws.Range("D34").Value = '''2.821'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.02%  '
$ws.Range("D35").Value = '''0.9589'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.27%  '
$ws.Range("E36").Value = '  +0.51%  '
$ws.Range("D37").Value = '''0.5781'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.19%  '
$ws.Range("D38").Value = '''0.01600'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.27%  '
$ws.Range("D39").Value = '''0.8626'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.67%  '
$ws.Range("D40").Value = '''5.833'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.46%  '
$ws.Range("E41").Value = '  +0.67%  '
$ws.Range("D42").Value = '1.043.93'
$ws.Range("E42").Value = '  +1.60%  '
$ws.Range("D43").Value = '''104.11'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.77%  '
$ws.Range("D44").Value = '1.801.88'
$ws.Range("E44").Value = '  +0.08%  '
$ws.Range("D45").Value = '''57.37'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.18%  '
$ws.Range("D46").Value = '0.0₈105'
$ws.Range("E46").Value = '  -5.17%  '
$ws.Range("D47").Value = '''1.010'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.26%  '
$ws.Range("E48").Value = '  +2.10%  '
$ws.Range("D49").Value = '''7.961'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.78%  '
$ws.Range("D50").Value = '''0.05163'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = '''1.430'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.80%  '
